# Removed unneeded semi colons from Python Examples
# Slide 25 ("Basic Python"), shape "TextBox 8" holds the Python code sample.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(25)
$shape = $s.Shapes.Item(3)
$tr = $shape.TextFrame.TextRange

# --- Change 1 -----------------------------------------------------------
# Paragraph: "tempInterest = float(interest) / 12;"
# Remove the trailing semicolon after "/ 12".
$tr.Characters(73, 1).Text = ""

# --- Change 2 -----------------------------------------------------------
# Paragraph: "   (tempInterest / (1.0 - ((1.0 + tempInterest) ** -float(term))));"
# Remove the trailing semicolon after the final "))))".
$tr.Characters(166, 1).Text = ""

# Split the remaining run "() ** -float(term))))" so that "term))))" becomes
# its own run (matching the authored edit, which left the new run's
# character run separate from the rest of the expression).
$newPart = $tr.Characters(158, 8)
$newPart.Font.Bold = $true
